## Applies the "Breadbox Poster" edit:
##   1. Moves "TextBox 33" (shape 9) down slightly   -> a:off/y  EMU 16022645 -> 16226329
##   2. Grows  "TextBox 38" (shape 12) taller         -> a:ext/cy EMU 2630312  -> 4846303
##   3. Replaces the body copy of "TextBox 38" with new placeholder text

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> point conversion (1 pt = 12700 EMU); PowerPoint's Shape geometry
# properties (Top/Left/Width/Height) are expressed in points.
$emuPerPt = 12700.0

# --- 1. "TextBox 33": nudge downward --------------------------------------
$textBox33 = $s.Shapes.Item(9)
$textBox33.Top = 16226329 / $emuPerPt

# --- 2 & 3. "TextBox 38": new copy + taller box ----------------------------
$textBox38 = $s.Shapes.Item(12)
$tr = $textBox38.TextFrame.TextRange

# The new paragraph, pre-split into the same run chunks as the authored
# slide (word-by-word, punctuation kept with its neighbour) so the
# resulting <a:r> run layout mirrors the original edit.
$runTexts = @(
    'As ',
    'verepe',
    ' ',
    'sectis',
    ' ',
    'enda',
    ' di ',
    'incto',
    ' ',
    'etusda',
    ' ',
    'delectet',
    ' ',
    'eaque',
    ' ',
    'simusda',
    ' ',
    'ndanimus',
    ', sin ',
    'cori',
    ' rem id ',
    'maximillut',
    ' ',
    'hari',
    ' con ',
    'ea',
    ' ',
    'quiam',
    ' ',
    'quam',
    ' et ',
    'enduci',
    ' ',
    'idunt',
    ', ',
    'corecta',
    ' ',
    'alicimp',
    ' ',
    'oriaeperum',
    ' ',
    'iusaeped',
    ' ',
    'unt',
    ' ',
    'quis',
    ' et ',
    'ut',
    ' ',
    'abor',
    ' ',
    'ame',
    ' ',
    'alis',
    ' ',
    'utatur',
    ' as senet a ',
    'quibus',
    ' ',
    'dolutatur',
    ' ',
    'sedit',
    ' ',
    'aliam',
    ' ',
    'apienda',
    ' ',
    'versper',
    ' ',
    'feriti',
    ' ',
    'delia',
    ' pa ',
    'dolo',
    ' ',
    'quamus',
    ' arum es ',
    'endellu',
    ' ',
    'pisquia',
    ' ',
    'aditatur',
    ' ',
    'audam',
    ' ',
    'earupid',
    ' quo ',
    'quunt',
    ' ',
    'odigniscil',
    ' ',
    'earum',
    ' ',
    'verios',
    ' ',
    'nonsequam',
    ' ',
    'duciis',
    ' et ',
    'apienda',
    ' ',
    'versper',
    ' ',
    'feriti',
    ' ',
    'delia',
    ' pa ',
    'dolo',
    ' ',
    'quamus',
    ' arum es ',
    'endellu',
    ' ',
    'pisquia',
    ' ',
    'aditatur',
    ' ',
    'audam',
    ' ',
    'earupid',
    ' quo ',
    'quunt',
    ' ',
    'odigniscil',
    ' ',
    'earum',
    ' ',
    'verios',
    ' ',
    'nonsequam'
)

$tr.Text = $runTexts[0]
for ($i = 1; $i -lt $runTexts.Count; $i++) {
    $null = $tr.InsertAfter($runTexts[$i])
}

# Set the height last - this textbox auto-fits to its text (a:spAutoFit),
# so the explicit Height assignment after the text edit is what pins the
# final box size instead of being overwritten by the auto-fit pass.
$textBox38.Height = 4846303 / $emuPerPt
